# Hortaliza, Vega Monumental Concepción - Coliflor
# Weekly update: insert the latest week's two new data rows (Primera / Segunda)
# at the top of the data block (row 41-42), pushing the rest of the
# historical rows down by two (old row 41 -> row 43, ..., old row 162 -> row 164).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 41; everything from row 41 down
# (through the old last row, 162) shifts down to rows 43-164.
$ws.Rows("41:42").Insert()

# --- New row 41: Coliflor, Primera, Región Metropolitana ---
$ws.Range("A41").Value = 11
$ws.Range("B41").Value = 'Vega Monumental Concepción'
$ws.Range("C41").Value = 'Bíobío'
$ws.Range("D41").Value = 44497
$ws.Range("E41").Value = 8
$ws.Range("F41").Value = 100112008
$ws.Range("G41").Value = 'Coliflor'
$ws.Range("H41").Value = 'Sin especificar'
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 2100
$ws.Range("K41").Value = 650
$ws.Range("L41").Value = 700
$ws.Range("M41").Value = 679
$ws.Range("N41").Value = '$/unidad'
$ws.Range("O41").Value = 'Región Metropolitana'
$ws.Range("P41").Value = 679
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = 'Hortaliza'

# --- New row 42: Coliflor, Segunda, Región Metropolitana ---
$ws.Range("A42").Value = 11
$ws.Range("B42").Value = 'Vega Monumental Concepción'
$ws.Range("C42").Value = 'Bíobío'
$ws.Range("D42").Value = 44497
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = 100112008
$ws.Range("G42").Value = 'Coliflor'
$ws.Range("H42").Value = 'Sin especificar'
$ws.Range("I42").Value = 'Segunda'
$ws.Range("J42").Value = 1200
$ws.Range("K42").Value = 600
$ws.Range("L42").Value = 600
$ws.Range("M42").Value = 600
$ws.Range("N42").Value = '$/unidad'
$ws.Range("O42").Value = 'Región Metropolitana'
$ws.Range("P42").Value = 600
$ws.Range("Q42").Value = 1
$ws.Range("R42").Value = 'Hortaliza'

Write-Output ("Dimension now: " + [string]$ws.Range("A1").Value)
